# "Generate Report for Handback" - localization-status.xlsx update
#
# The CI run that hands a translation back now records:
#   - the overall status text ("Handed back: in sync with en-US" instead of
#     "Ready for handoff") on the Overview sheet and on each language sheet,
#   - the resolved "Latest Target File" (a link back to the source .md),
#   - the "Latest Handback File" (the xlf that was generated), and
#   - the "Latest Handback DateTime" timestamp
# for both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$sourceDisplay = "a450824d-a77d-4865-be7d-56a2dd79518d.md"
$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f03163124286c95cd923a1d82bd47b8eef9460c1/e2e/a450824d-a77d-4865-be7d-56a2dd79518d.md"

# ---- Overview sheet: refresh the per-language status cells ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText

# Latest Target File (I2) now links back to the source document.
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $sourceUrl, "", "", $sourceDisplay)

# Latest Handback File (J2) records the generated xliff file name.
$wsZh.Range("J2").Value = "a450824d-a77d-4865-be7d-56a2dd79518d.041de7cc909284d1a1c6d39e276c33ca8ea8a71a.zh-cn.xlf"

# Latest Handback DateTime (K2).
$wsZh.Range("K2").Value = "2016-08-18 06:57:41"

$wsZh.Columns.Item(3).ColumnWidth = 29.2
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText

# Latest Target File (I2) now links back to the source document.
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $sourceUrl, "", "", $sourceDisplay)

# Latest Handback File (J2) records the generated xliff file name.
$wsDe.Range("J2").Value = "a450824d-a77d-4865-be7d-56a2dd79518d.041de7cc909284d1a1c6d39e276c33ca8ea8a71a.de-de.xlf"

# Latest Handback DateTime (K2).
$wsDe.Range("K2").Value = "2016-08-18 06:57:48"

$wsDe.Columns.Item(3).ColumnWidth = 29.2
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17

# ---- Overview sheet: widen the now-longer status columns ----
$wsOverview.Columns.Item(5).ColumnWidth = 29.2
$wsOverview.Columns.Item(6).ColumnWidth = 29.2
